$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must remain text (matches source formatting)
$textCells = @("D5", "D6", "D8", "D11", "D12", "D13", "D14", "D16", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D30", "D31", "D32", "D33", "D39", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range('D2').Value = '89.695.08'
$ws.Range('E2').Value = '  +0.91%  '
$ws.Range('D3').Value = '3.047.87'
$ws.Range('E3').Value = '  -2.95%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '211.30'
$ws.Range('E5').Value = '  -1.78%  '
$ws.Range('D6').Value = '612.30'
$ws.Range('E6').Value = '  -3.37%  '
$ws.Range('E7').Value = '  -8.18%  '
$ws.Range('D8').Value = '0.884'
$ws.Range('E8').Value = '  +18.32%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').Value = '3.047.82'
$ws.Range('E10').Value = '  -2.90%  '
$ws.Range('D11').Value = '0.676'
$ws.Range('E11').Value = '  +21.44%  '
$ws.Range('D12').Value = '0.188'
$ws.Range('E12').Value = '  +5.01%  '
$ws.Range('D13').Value = '0.0000239'
$ws.Range('E13').Value = '  -4.80%  '
$ws.Range('D14').Value = '5.37'
$ws.Range('E14').Value = '  +1.16%  '
$ws.Range('D15').Value = '89.564.32'
$ws.Range('E15').Value = '  +0.90%  '
$ws.Range('D16').Value = '32.08'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').Value = '3.590.80'
$ws.Range('E17').Value = '  -3.54%  '
$ws.Range('D18').Value = '3.045.35'
$ws.Range('E18').Value = '  -3.38%  '
$ws.Range('D19').Value = '3.32'
$ws.Range('E19').Value = '  -2.36%  '
$ws.Range('D20').Value = '0.0000216'
$ws.Range('E20').Value = '  -5.91%  '
$ws.Range('D21').Value = '13.41'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('D22').Value = '425.07'
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '8.21'
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('B24').Value = 'Polkadot'
$ws.Range('C24').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D24').Value = '5.01'
$ws.Range('E24').Value = '  +2.35%  '
$ws.Range('D25').Value = '5.38'
$ws.Range('E25').Value = '  -0.91%  '
$ws.Range('D26').Value = '83.92'
$ws.Range('E26').Value = '  +2.92%  '
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D31').Value = '8.20'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = '504.39'
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('D33').Value = '3.71'
$ws.Range('E33').Value = '  -8.44%  '
$ws.Range('E34').Value = '  -6.76%  '
$ws.Range('E35').Value = '  +4.31%  '
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('E37').Value = '  -5.13%  '
$ws.Range('E38').Value = '  -9.24%  '
$ws.Range('D39').Value = '22.26'
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  -0.82%  '
$ws.Range('E43').Value = '  +7.11%  '
$ws.Range('E44').Value = '  -2.01%  '
$ws.Range('D45').Value = '147.08'
$ws.Range('E45').Value = '  +0.82%  '
$ws.Range('D46').Value = '43.29'
$ws.Range('E46').Value = '  -1.08%  '
$ws.Range('D47').Value = '0.0687'
$ws.Range('E47').Value = '  +11.26%  '
$ws.Range('D48').Value = '4.09'
$ws.Range('E48').Value = '  +4.05%  '
$ws.Range('B49').Value = 'Aave'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D49').Value = '160.76'
$ws.Range('E49').Value = '  -2.66%  '
$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = '1.21'
$ws.Range('E50').Value = '  +2.12%  '
$ws.Range('D51').Value = '0.697'
$ws.Range('E51').Value = '  -3.38%  '
